$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 40; this shifts existing rows 40-165 down to 41-166
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with its values
$ws.Cells.Item(40, 1).Value = 5
$ws.Cells.Item(40, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(40, 3).Value = "Maule"
$ws.Cells.Item(40, 4).Value = 44453
$ws.Cells.Item(40, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 5).Value = 7
$ws.Cells.Item(40, 6).Value = 100112006
$ws.Cells.Item(40, 7).Value = "Repollo"
$ws.Cells.Item(40, 8).Value = "Crespo record"
$ws.Cells.Item(40, 9).Value = "Segunda"
$ws.Cells.Item(40, 10).Value = 5000
$ws.Cells.Item(40, 11).Value = 300
$ws.Cells.Item(40, 12).Value = 300
$ws.Cells.Item(40, 13).Value = 300
$ws.Cells.Item(40, 14).Value = "$/unidad"
$ws.Cells.Item(40, 15).Value = "Región del Maule"
$ws.Cells.Item(40, 16).Value = 300
$ws.Cells.Item(40, 17).Value = 1
$ws.Cells.Item(40, 18).Value = "Hortaliza"
